$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.323.47'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '1.861.83'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4768'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2749'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06448'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").Value = '1.897.91'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07425'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.001'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '86.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6328'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("D16").Value = '30.303.15'
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9994'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007396'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.05%  '
$ws.Range("D21").Value = '2.097.10'
$ws.Range("E21").Value = '  -6.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9988'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.114'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.3914'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.011'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.308'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  -3.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.861'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.60%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1013'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.58%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.383'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.239'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.918'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.153'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7266'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9988'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.690'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01945'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.631'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9097'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.992'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '105.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4126'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.567'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.080'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '61.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.05%  '
$ws.Range("E49").Value = '  -5.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.764'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.406'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.95%  '
